$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 123.75
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 180
$ws.Range("K2").Value = 30
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = 83
$ws.Range("N2").Value = -406

$ws.Range("H132").Value = 779646.9399999999
$ws.Range("I132").Value = 1689.9474
$ws.Range("J132").Value = 8170238.5
$ws.Range("K132").Value = 5069.8422
$ws.Range("L132").Value = 24510715.5
$ws.Range("M132").Value = -2539.8422
$ws.Range("N132").Value = -24515775.5

$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200

$ws.Range("H138").Value = 2316419
$ws.Range("I138").Value = 1168.7931
$ws.Range("J138").Value = 11908170
$ws.Range("K138").Value = 3506.379300000001
$ws.Range("L138").Value = 35724510
$ws.Range("M138").Value = 1633.620699999999
$ws.Range("N138").Value = -35734790

$ws.Range("H141").Value = 2993.25
$ws.Range("I141").Value = 2068.8076
$ws.Range("J141").Value = 6999.1665
$ws.Range("K141").Value = 6206.4228
$ws.Range("L141").Value = 20997.4995
$ws.Range("M141").Value = -1026.4228
$ws.Range("N141").Value = -31357.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5033.3335
$ws.Range("I2").Value = 8000
$ws.Range("J2").Value = 4440
$ws.Range("K2").Value = 8000
$ws.Range("L2").Value = 4440
$ws.Range("M2").Value = -7887
$ws.Range("N2").Value = -4666

$ws.Range("H32").Value = 2199.84
$ws.Range("I32").Value = 1644.3457
$ws.Range("J32").Value = 4568
$ws.Range("K32").Value = 1644.3457
$ws.Range("L32").Value = 4568
$ws.Range("M32").Value = -1357.3457
$ws.Range("N32").Value = -5142

$ws.Range("H47").Value = 21999.334
$ws.Range("J47").Value = 21999.334
$ws.Range("L47").Value = 21999.334
$ws.Range("N47").Value = -23449.334

$ws.Range("H61").Value = 21784022
$ws.Range("I61").Value = 27055312
$ws.Range("K61").Value = 27055312
$ws.Range("M61").Value = -27055100

$ws.Range("H74").Value = 6576321.5
$ws.Range("I74").Value = 8155403.5
$ws.Range("J74").Value = 102084.4
$ws.Range("K74").Value = 8155403.5
$ws.Range("L74").Value = 102084.4
$ws.Range("M74").Value = -8154529.5
$ws.Range("N74").Value = -103832.4

$ws.Range("H77").Value = 6576321.5
$ws.Range("I77").Value = 8155403.5
$ws.Range("J77").Value = 102084.4
$ws.Range("K77").Value = 40777017.5
$ws.Range("L77").Value = 510422
$ws.Range("M77").Value = -40772649.5
$ws.Range("N77").Value = -519158

$ws.Range("H116").Value = 5033.3335
$ws.Range("I116").Value = 8000
$ws.Range("J116").Value = 4440
$ws.Range("K116").Value = 8000
$ws.Range("L116").Value = 4440
$ws.Range("M116").Value = -5706
$ws.Range("N116").Value = -9028

$ws.Range("H132").Value = 35151.984
$ws.Range("I132").Value = 25653.85
$ws.Range("J132").Value = 53243.668
$ws.Range("K132").Value = 76961.54999999999
$ws.Range("L132").Value = 159731.004
$ws.Range("M132").Value = -74431.54999999999
$ws.Range("N132").Value = -164791.004

$ws.Range("H136").Value = 21784022
$ws.Range("I136").Value = 27055312
$ws.Range("K136").Value = 81165936
$ws.Range("M136").Value = -81163386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5033.3335
$ws.Range("I3").Value = 8000
$ws.Range("J3").Value = 4440
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 4440
$ws.Range("M3").Value = -7886
$ws.Range("N3").Value = -4668

$ws.Range("H9").Value = 30000
$ws.Range("J9").Value = 30000
$ws.Range("L9").Value = 30000
$ws.Range("N9").Value = -30336

$ws.Range("H44").Value = 15000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 111112456
$ws.Range("J16").Value = 166667680
$ws.Range("L16").Value = 166667680
$ws.Range("N16").Value = -166668254

$ws.Range("H31").Value = 3412
$ws.Range("I31").Value = 1477.8572
$ws.Range("J31").Value = 6313.2144
$ws.Range("K31").Value = 1477.8572
$ws.Range("L31").Value = 6313.2144
$ws.Range("M31").Value = -1182.8572
$ws.Range("N31").Value = -6903.2144

$ws.Range("H34").Value = 3412
$ws.Range("I34").Value = 1477.8572
$ws.Range("J34").Value = 6313.2144
$ws.Range("K34").Value = 1477.8572
$ws.Range("L34").Value = 6313.2144
$ws.Range("M34").Value = -1275.8572
$ws.Range("N34").Value = -6717.2144

$ws.Range("H58").Value = 17858620
$ws.Range("I58").Value = 22728566
$ws.Range("J58").Value = 2149.6667
$ws.Range("K58").Value = 22728566
$ws.Range("L58").Value = 2149.6667
$ws.Range("M58").Value = -22728363
$ws.Range("N58").Value = -2555.6667

$ws.Range("H113").Value = 111112456
$ws.Range("J113").Value = 166667680
$ws.Range("L113").Value = 166667680
$ws.Range("N113").Value = -166672020

$ws.Range("H132").Value = 15031.486
$ws.Range("I132").Value = 994.322
$ws.Range("J132").Value = 78738.62
$ws.Range("K132").Value = 2982.966
$ws.Range("L132").Value = 236215.86
$ws.Range("M132").Value = -452.9659999999999
$ws.Range("N132").Value = -241275.86

$ws.Range("H134").Value = 16723.783
$ws.Range("I134").Value = 1028.58
$ws.Range("J134").Value = 58026.95
$ws.Range("K134").Value = 3085.74
$ws.Range("L134").Value = 174080.85
$ws.Range("M134").Value = -550.7399999999998
$ws.Range("N134").Value = -179150.85

$ws.Range("H136").Value = 17858620
$ws.Range("I136").Value = 22728566
$ws.Range("J136").Value = 2149.6667
$ws.Range("K136").Value = 68185698
$ws.Range("L136").Value = 6449.000100000001
$ws.Range("M136").Value = -68183148
$ws.Range("N136").Value = -11549.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1000.19116
$ws.Range("I131").Value = 416.25
$ws.Range("J131").Value = 1078.05
$ws.Range("K131").Value = 1248.75
$ws.Range("L131").Value = 3234.15
$ws.Range("M131").Value = 3791.25
$ws.Range("N131").Value = -13314.15

$ws.Range("H132").Value = 2824.2104
$ws.Range("I132").Value = 2091.6667
$ws.Range("K132").Value = 18825.0003
$ws.Range("M132").Value = -16295.0003

$ws.Range("H137").Value = 36650.53
$ws.Range("I137").Value = 764.75
$ws.Range("K137").Value = 2294.25
$ws.Range("M137").Value = 2805.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 15000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 15000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 15000
$ws.Range("M44").Value = ""
$ws.Range("N44").Value = -16192

$ws.Range("H126").Value = 1955.4445
$ws.Range("I126").Value = 1599.8
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 4799.4
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -2329.4
$ws.Range("N126").Value = -12140

$ws.Range("H132").Value = 36260.81
$ws.Range("I132").Value = 24197.256
$ws.Range("J132").Value = 73313.14
$ws.Range("K132").Value = 72591.76800000001
$ws.Range("L132").Value = 219939.42
$ws.Range("M132").Value = -70061.76800000001
$ws.Range("N132").Value = -224999.42

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1697.8846
$ws.Range("I100").Value = 1624.5
$ws.Range("K100").Value = 1624.5
$ws.Range("M100").Value = -1083.5

$ws.Range("H136").Value = 41399.12
$ws.Range("I136").Value = 27422.553
$ws.Range("K136").Value = 82267.659
$ws.Range("M136").Value = -79717.659

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 67131.09
$ws.Range("I132").Value = 59304.05
$ws.Range("J132").Value = 78570.62
$ws.Range("K132").Value = 177912.15
$ws.Range("L132").Value = 235711.86
$ws.Range("M132").Value = -175382.15
$ws.Range("N132").Value = -240771.86

$ws.Range("H136").Value = 23081.076
$ws.Range("I136").Value = 16778.645
$ws.Range("J136").Value = 36106.1
$ws.Range("K136").Value = 50335.935
$ws.Range("L136").Value = 108318.3
$ws.Range("M136").Value = -47785.935
$ws.Range("N136").Value = -113418.3
